# ------------------------------------------------------------------
# PlayerPerformance_4584.xlsx edit:
#  - Add "Player Info" sheet (before "ODI Batting")
#  - "ODI Batting": rename MATCH_CARD_LINK -> MATCH_CODE, strip URL to
#    bare match code, drop stray empty INNING_NUMBER cells
#  - "ODI Bowling": rename MATCH_CARD_LINK -> MATCH_CODE, strip URL to
#    bare match code
#  - Add "ODI Batting Extra" sheet (after "ODI Bowling")
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

function Test-LooksNumeric([string]$s) {
    # Strings that Excel would otherwise auto-coerce into a number,
    # percentage or date when assigned through .Value.
    if ($s -match '^-?\d+(\.\d+)?%?$') { return $true }
    if ($s -match '^\d{1,2}/\d{1,2}/\d{4}$') { return $true }
    return $false
}

function Set-TextValue($range, [string]$text) {
    # Write a value that must stay literal text. Only prefix with an
    # apostrophe (forcing text storage) when the text would otherwise
    # be auto-coerced into a number/date/percent by Excel.
    if (Test-LooksNumeric $text) {
        $range.Value = "'" + $text
    } else {
        $range.Value = $text
    }
}

# ------------------------------------------------------------------
# 1) ODI Batting : MATCH_CARD_LINK -> MATCH_CODE
# ------------------------------------------------------------------
$wsBatting = $wb.Worksheets.Item("ODI Batting")
Set-TextValue $wsBatting.Range("D1") "MATCH_CODE"

$battingRows = 2..26
foreach ($r in $battingRows) {
    $cell = $wsBatting.Cells.Item($r, 4)
    $link = [string]$cell.Value2
    if ($link -and $link -ne "") {
        $parts = $link -split "MatchCode="
        if ($parts.Length -gt 1) {
            Set-TextValue $cell $parts[1]
        }
    }
}

# Rows whose INNING_NUMBER (col B) is blank carried a stray empty cell
# in the original file; clear them so they disappear entirely.
$blankInningRows = @(4,5,6,7,10,12,13,16,17,19,20,21,24,25,26)
foreach ($r in $blankInningRows) {
    $wsBatting.Cells.Item($r, 2).Value = ""
}

# ------------------------------------------------------------------
# 2) ODI Bowling : MATCH_CARD_LINK -> MATCH_CODE
# ------------------------------------------------------------------
$wsBowling = $wb.Worksheets.Item("ODI Bowling")
Set-TextValue $wsBowling.Range("B1") "MATCH_CODE"

$bowlingRows = 2..24
foreach ($r in $bowlingRows) {
    $cell = $wsBowling.Cells.Item($r, 2)
    $link = [string]$cell.Value2
    if ($link -and $link -ne "") {
        $parts = $link -split "MatchCode="
        if ($parts.Length -gt 1) {
            Set-TextValue $cell $parts[1]
        }
    }
}

# ------------------------------------------------------------------
# 3) New sheet "Player Info" (inserted before "ODI Batting")
# ------------------------------------------------------------------
$wsBatting = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($wsBatting)
$playerInfo.Name = "Player Info"

$piHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $piHeaders.Length; $c++) {
    $cell = $playerInfo.Cells.Item(1, $c)
    Set-TextValue $cell $piHeaders[$c - 1]
}
$piHeaderRange = $playerInfo.Range("A1:D1")
$piHeaderRange.Font.Bold = $true
$piHeaderRange.HorizontalAlignment = -4108
$piHeaderRange.VerticalAlignment = -4160
$piHeaderRange.Borders.LineStyle = 1

$piRow = @("4584", "Joshua Brian Little", "Right Handed", "Left Arm Fast")
for ($c = 1; $c -le $piRow.Length; $c++) {
    $cell = $playerInfo.Cells.Item(2, $c)
    Set-TextValue $cell $piRow[$c - 1]
}

# ------------------------------------------------------------------
# 4) New sheet "ODI Batting Extra" (appended after "ODI Bowling")
# ------------------------------------------------------------------
$wsBowling = $wb.Worksheets.Item("ODI Bowling")
$extra = $wb.Worksheets.Add($null, $wsBowling)
$extra.Name = "ODI Batting Extra"

$exHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $exHeaders.Length; $c++) {
    $cell = $extra.Cells.Item(1, $c)
    Set-TextValue $cell $exHeaders[$c - 1]
}
$exHeaderRange = $extra.Range("A1:F1")
$exHeaderRange.Font.Bold = $true
$exHeaderRange.HorizontalAlignment = -4108
$exHeaderRange.VerticalAlignment = -4160
$exHeaderRange.Borders.LineStyle = 1

# MATCH_CODE, BATTING_POSITION (numeric or blank), NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$exData = @(
    @("4428", 10,   "1", "0", "4.64%", "NO"),
    @("4448", 10,   "0", "0", "0.61%", "NO"),
    @("4466", $null, $null, $null, $null, "NO"),
    @("4467", 10,   $null, $null, $null, "NO"),
    @("4468", $null, $null, $null, $null, "NO"),
    @("4474", $null, $null, $null, $null, "NO"),
    @("4475", $null, $null, $null, $null, "NO"),
    @("4478", $null, $null, $null, $null, "NO"),
    @("4492", 11,   $null, $null, $null, "NO"),
    @("4494", 11,   "0", "0", $null, "NO"),
    @("4496", $null, $null, $null, $null, "NO"),
    @("4519", 11,   $null, $null, $null, "YES"),
    @("4520", $null, $null, $null, $null, "NO"),
    @("4522", $null, $null, $null, $null, $null),
    @("4605", $null, $null, $null, $null, $null),
    @("4608", $null, $null, $null, $null, $null),
    @("4614", $null, $null, $null, $null, $null),
    @("4693", $null, $null, $null, $null, $null),
    @("4694", $null, $null, $null, $null, $null),
    @("4696", $null, $null, $null, $null, $null)
)

$r = 2
foreach ($row in $exData) {
    Set-TextValue $extra.Cells.Item($r, 1) $row[0]

    $posCell = $extra.Cells.Item($r, 2)
    if ($null -eq $row[1]) { $posCell.Value = "" } else { $posCell.Value = $row[1] }

    for ($c = 3; $c -le 5; $c++) {
        $val = $row[$c - 1]
        $cell = $extra.Cells.Item($r, $c)
        if ($null -eq $val) { $cell.Value = "" } else { Set-TextValue $cell $val }
    }

    $momCell = $extra.Cells.Item($r, 6)
    if ($null -eq $row[5]) { $momCell.Value = "" } else { Set-TextValue $momCell $row[5] }

    $r++
}
